$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry represents the updated Coin/Link/Price/Volume(1h) values for a data row,
# reflecting a refreshed snapshot of the cryptos list.
$updates = @(
    @{ Row = 2; B = $null; C = $null; D = '29.240.39'; E = '  +0.03%  ' },
    @{ Row = 3; B = $null; C = $null; D = '1.856.78'; E = '  -0.29%  ' },
    @{ Row = 4; B = $null; C = $null; D = '0.9991'; E = '  -0.22%  ' },
    @{ Row = 5; B = $null; C = $null; D = '0.7016'; E = '  -0.59%  ' },
    @{ Row = 6; B = $null; C = $null; D = '241.40'; E = '  -0.39%  ' },
    @{ Row = 7; B = $null; C = $null; D = '0.9995'; E = '  -0.22%  ' },
    @{ Row = 8; B = 'Dogecoin'; C = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D = '0.07751'; E = '  -0.85%  ' },
    @{ Row = 9; B = 'Cardano'; C = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D = '0.3094'; E = '  -0.54%  ' },
    @{ Row = 10; B = $null; C = $null; D = '23.85'; E = '  -1.75%  ' },
    @{ Row = 11; B = $null; C = $null; D = '0.07798'; E = '  -2.51%  ' },
    @{ Row = 12; B = $null; C = $null; D = '1.867.37'; E = '  -2.19%  ' },
    @{ Row = 13; B = $null; C = $null; D = '5.108'; E = '  -1.32%  ' },
    @{ Row = 14; B = $null; C = $null; D = '92.28'; E = '  -1.41%  ' },
    @{ Row = 15; B = $null; C = $null; D = '0.6889'; E = '  -0.87%  ' },
    @{ Row = 16; B = $null; C = $null; D = '6.513'; E = '  +2.62%  ' },
    @{ Row = 17; B = $null; C = $null; D = '0.000008428'; E = '  +1.93%  ' },
    @{ Row = 18; B = $null; C = $null; D = '29.246.34'; E = '  -0.34%  ' },
    @{ Row = 19; B = $null; C = $null; D = '249.88'; E = '  -0.71%  ' },
    @{ Row = 20; B = $null; C = $null; D = '2.115.50'; E = '  -2.63%  ' },
    @{ Row = 21; B = $null; C = $null; D = '12.88'; E = '  -1.70%  ' },
    @{ Row = 22; B = $null; C = $null; D = $null; E = '  -0.17%  ' },
    @{ Row = 23; B = $null; C = $null; D = '7.535'; E = '  -0.11%  ' },
    @{ Row = 24; B = $null; C = $null; D = '0.9992'; E = '  -0.25%  ' },
    @{ Row = 25; B = $null; C = $null; D = '0.1518'; E = '  -2.36%  ' },
    @{ Row = 26; B = $null; C = $null; D = '160.41'; E = '  +0.57%  ' },
    @{ Row = 27; B = $null; C = $null; D = '8.862'; E = '  -1.42%  ' },
    @{ Row = 28; B = $null; C = $null; D = '18.53'; E = '  -0.93%  ' },
    @{ Row = 29; B = $null; C = $null; D = '1.562'; E = '  +4.30%  ' },
    @{ Row = 30; B = $null; C = $null; D = '4.244'; E = '  -0.58%  ' },
    @{ Row = 31; B = $null; C = $null; D = '4.198'; E = '  -1.66%  ' },
    @{ Row = 32; B = $null; C = $null; D = '1.195'; E = '  -1.35%  ' },
    @{ Row = 33; B = $null; C = $null; D = '0.05200'; E = '  -1.29%  ' },
    @{ Row = 34; B = $null; C = $null; D = '0.7646'; E = '  +2.66%  ' },
    @{ Row = 35; B = $null; C = $null; D = '1.848'; E = '  -2.14%  ' },
    @{ Row = 36; B = $null; C = $null; D = '1.165'; E = '  +0.73%  ' },
    @{ Row = 37; B = $null; C = $null; D = '2.710'; E = '  +0.14%  ' },
    @{ Row = 38; B = $null; C = $null; D = '0.01864'; E = '  +0.31%  ' },
    @{ Row = 39; B = $null; C = $null; D = '1.220.45'; E = '  -2.38%  ' },
    @{ Row = 40; B = $null; C = $null; D = '2.723'; E = '  -0.61%  ' },
    @{ Row = 41; B = $null; C = $null; D = '0.8968'; E = '  -0.60%  ' },
    @{ Row = 42; B = $null; C = $null; D = '109.90'; E = '  -1.16%  ' },
    @{ Row = 43; B = $null; C = $null; D = '0.9988'; E = '  -0.21%  ' },
    @{ Row = 44; B = $null; C = $null; D = '5.550'; E = '  -11.83%  ' },
    @{ Row = 45; B = $null; C = $null; D = '2.012.93'; E = '  -2.75%  ' },
    @{ Row = 46; B = 'Aave'; C = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D = '65.38'; E = '  -9.07%  ' },
    @{ Row = 47; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '9.553'; E = '  +1.84%  ' },
    @{ Row = 48; B = $null; C = $null; D = '0.5177'; E = '  -0.46%  ' },
    @{ Row = 49; B = 'BabyDogeCoin'; C = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; D = '0.00000000122'; E = '  -5.35%  ' },
    @{ Row = 50; B = $null; C = $null; D = '1.755'; E = '  -1.66%  ' },
    @{ Row = 51; B = $null; C = $null; D = '7.034'; E = '  +0.71%  ' }
)

function Looks-Numeric([string]$s) {
    return ($s -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$')
}

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.B -ne $null) {
        $ws.Cells.Item($row, 2).Value = $u.B
    }
    if ($u.C -ne $null) {
        $ws.Cells.Item($row, 3).Value = $u.C
    }
    if ($u.D -ne $null) {
        $cell = $ws.Cells.Item($row, 4)
        if (Looks-Numeric $u.D) {
            # Force text storage so values like "0.9991" or "0.000008428" keep their
            # exact original textual representation instead of becoming a Double.
            $cell.NumberFormat = "@"
            $cell.Value = $u.D
            $cell.ClearFormats()
        } else {
            $cell.Value = $u.D
        }
    }
    if ($u.E -ne $null) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
